{"js": "// Replace each arithmetic expression in the worksheet table with its\n// updated version. Each old expression text is unique within the document,\n// so a targeted search + replace for each pair reproduces the diff exactly\n// while leaving every other paragraph/run/formatting untouched.\nconst pairs = [[\"67-8=\", \"67+9=\"], [\"33-0=\", \"33+26=\"], [\"93-40=\", \"91-63=\"], [\"42+17=\", \"24+46=\"], [\"35+38=\", \"94-9=\"], [\"76-33=\", \"67+28=\"], [\"51-24=\", \"21+32=\"], [\"53-5=\", \"60-27=\"], [\"98-7=\", \"82-52=\"], [\"45-12=\", \"48-43=\"], [\"72-67=\", \"16+12=\"], [\"96-5=\", \"72+13=\"], [\"64-9=\", \"14+48=\"], [\"10+30=\", \"20+40=\"], [\"76-49=\", \"41-26=\"], [\"42+35=\", \"98-28=\"], [\"6+57=\", \"76-29=\"], [\"90-8=\", \"39+10=\"], [\"35+44=\", \"49-46=\"], [\"58+32=\", \"47-44=\"], [\"2+30=\", \"95-44=\"], [\"29+29=\", \"57-55=\"], [\"71+14=\", \"10+21=\"], [\"99-27=\", \"81+16=\"], [\"2+2=\", \"29+34=\"], [\"75-27=\", \"81-40=\"], [\"91-84=\", \"21+74=\"], [\"43-7=\", \"16-2=\"], [\"9+2=\", \"43-14=\"], [\"85-69=\", \"9-5=\"], [\"79-39=\", \"6+28=\"], [\"5+74=\", \"43-12=\"], [\"80+13=\", \"96-3=\"], [\"72+23=\", \"9+5=\"], [\"45+42=\", \"79-54=\"], [\"81-36=\", \"43+47=\"], [\"98-45=\", \"2+71=\"], [\"15+60=\", \"92-11=\"], [\"51+26=\", \"0+58=\"], [\"96-40=\", \"17+12=\"], [\"0+34=\", \"22+41=\"], [\"8+75=\", \"21+50=\"], [\"29+48=\", \"77-25=\"], [\"60+2=\", \"48-34=\"], [\"73-70=\", \"66-26=\"], [\"13+7=\", \"68-24=\"], [\"91-1=\", \"15-14=\"], [\"6+49=\", \"45-30=\"], [\"21+11=\", \"20+49=\"], [\"24+31=\", \"84-18=\"], [\"44+23=\", \"20+17=\"], [\"61-1=\", \"33+59=\"], [\"45+50=\", \"31-25=\"], [\"81-27=\", \"59-9=\"], [\"50-46=\", \"81-12=\"], [\"86-26=\", \"8+60=\"], [\"17+36=\", \"24+30=\"], [\"70+2=\", \"17+73=\"], [\"79-38=\", \"47+19=\"], [\"31+62=\", \"37+5=\"], [\"23+2=\", \"92-28=\"], [\"69-27=\", \"87+11=\"], [\"80-20=\", \"90-47=\"], [\"67-56=\", \"64-50=\"], [\"76-51=\", \"12+85=\"], [\"16+33=\", \"16+77=\"], [\"91-51=\", \"60+22=\"], [\"29+24=\", \"24+58=\"], [\"62+6=\", \"42-20=\"], [\"84-36=\", \"46-38=\"], [\"2+46=\", \"86-74=\"], [\"10+76=\", \"36-4=\"], [\"63-2=\", \"95-52=\"], [\"53-20=\", \"24+67=\"], [\"3+60=\", \"33+1=\"], [\"89-7=\", \"55+43=\"], [\"32+38=\", \"50+24=\"], [\"37+51=\", \"72+17=\"], [\"95-45=\", \"85-71=\"], [\"18-11=\", \"11+68=\"], [\"11+32=\", \"72+27=\"], [\"31+38=\", \"69-4=\"], [\"36+9=\", \"5+14=\"], [\"52-24=\", \"84+14=\"], [\"89-70=\", \"97-30=\"], [\"58-31=\", \"80-11=\"], [\"30-11=\", \"0+53=\"], [\"21-7=\", \"90-61=\"], [\"66+8=\", \"13+28=\"], [\"49-31=\", \"97-37=\"], [\"91+2=\", \"22+33=\"], [\"49-13=\", \"31+55=\"], [\"8+16=\", \"25+61=\"], [\"34+14=\", \"51-8=\"], [\"69-43=\", \"62-39=\"], [\"73+21=\", \"13+55=\"], [\"77-49=\", \"30+0=\"], [\"95-83=\", \"63+32=\"], [\"18-12=\", \"1+85=\"], [\"60-55=\", \"39-28=\"]];\n\nconst body = context.document.body;\n\n// Phase 1: issue all searches, then resolve them with a single sync so we\n// don't pay a network/RPC round-trip per pair.\nconst searchResults = pairs.map(([oldText]) => {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  return results;\n});\n\nawait context.sync();\n\n// Phase 2: replace the matched range's text in place, preserving the run's\n// existing formatting (font, size, etc.) since insertText(..., \"Replace\")\n// only swaps the text content of the matched range.\nfor (let i = 0; i < pairs.length; i++) {\n  const [oldText, newText] = pairs[i];\n  const results = searchResults[i];\n  if (results.items.length === 0) {\n    throw new Error(`Could not find expected text: ${oldText}`);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update each worksheet cell's arithmetic expression to its revised\n# value. Every 'old' expression string is unique in the document, so a\n# plain Find/Replace (wdReplaceAll, but each pattern only ever matches\n# once) swaps exactly the target <w:t> run and leaves formatting intact.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('67-8=', '67+9=')\n    ,@('33-0=', '33+26=')\n    ,@('93-40=', '91-63=')\n    ,@('42+17=', '24+46=')\n    ,@('35+38=', '94-9=')\n    ,@('76-33=', '67+28=')\n    ,@('51-24=', '21+32=')\n    ,@('53-5=', '60-27=')\n    ,@('98-7=', '82-52=')\n    ,@('45-12=', '48-43=')\n    ,@('72-67=', '16+12=')\n    ,@('96-5=', '72+13=')\n    ,@('64-9=', '14+48=')\n    ,@('10+30=', '20+40=')\n    ,@('76-49=', '41-26=')\n    ,@('42+35=', '98-28=')\n    ,@('6+57=', '76-29=')\n    ,@('90-8=', '39+10=')\n    ,@('35+44=', '49-46=')\n    ,@('58+32=', '47-44=')\n    ,@('2+30=', '95-44=')\n    ,@('29+29=', '57-55=')\n    ,@('71+14=', '10+21=')\n    ,@('99-27=', '81+16=')\n    ,@('2+2=', '29+34=')\n    ,@('75-27=', '81-40=')\n    ,@('91-84=', '21+74=')\n    ,@('43-7=', '16-2=')\n    ,@('9+2=', '43-14=')\n    ,@('85-69=', '9-5=')\n    ,@('79-39=', '6+28=')\n    ,@('5+74=', '43-12=')\n    ,@('80+13=', '96-3=')\n    ,@('72+23=', '9+5=')\n    ,@('45+42=', '79-54=')\n    ,@('81-36=', '43+47=')\n    ,@('98-45=', '2+71=')\n    ,@('15+60=', '92-11=')\n    ,@('51+26=', '0+58=')\n    ,@('96-40=', '17+12=')\n    ,@('0+34=', '22+41=')\n    ,@('8+75=', '21+50=')\n    ,@('29+48=', '77-25=')\n    ,@('60+2=', '48-34=')\n    ,@('73-70=', '66-26=')\n    ,@('13+7=', '68-24=')\n    ,@('91-1=', '15-14=')\n    ,@('6+49=', '45-30=')\n    ,@('21+11=', '20+49=')\n    ,@('24+31=', '84-18=')\n    ,@('44+23=', '20+17=')\n    ,@('61-1=', '33+59=')\n    ,@('45+50=', '31-25=')\n    ,@('81-27=', '59-9=')\n    ,@('50-46=', '81-12=')\n    ,@('86-26=', '8+60=')\n    ,@('17+36=', '24+30=')\n    ,@('70+2=', '17+73=')\n    ,@('79-38=', '47+19=')\n    ,@('31+62=', '37+5=')\n    ,@('23+2=', '92-28=')\n    ,@('69-27=', '87+11=')\n    ,@('80-20=', '90-47=')\n    ,@('67-56=', '64-50=')\n    ,@('76-51=', '12+85=')\n    ,@('16+33=', '16+77=')\n    ,@('91-51=', '60+22=')\n    ,@('29+24=', '24+58=')\n    ,@('62+6=', '42-20=')\n    ,@('84-36=', '46-38=')\n    ,@('2+46=', '86-74=')\n    ,@('10+76=', '36-4=')\n    ,@('63-2=', '95-52=')\n    ,@('53-20=', '24+67=')\n    ,@('3+60=', '33+1=')\n    ,@('89-7=', '55+43=')\n    ,@('32+38=', '50+24=')\n    ,@('37+51=', '72+17=')\n    ,@('95-45=', '85-71=')\n    ,@('18-11=', '11+68=')\n    ,@('11+32=', '72+27=')\n    ,@('31+38=', '69-4=')\n    ,@('36+9=', '5+14=')\n    ,@('52-24=', '84+14=')\n    ,@('89-70=', '97-30=')\n    ,@('58-31=', '80-11=')\n    ,@('30-11=', '0+53=')\n    ,@('21-7=', '90-61=')\n    ,@('66+8=', '13+28=')\n    ,@('49-31=', '97-37=')\n    ,@('91+2=', '22+33=')\n    ,@('49-13=', '31+55=')\n    ,@('8+16=', '25+61=')\n    ,@('34+14=', '51-8=')\n    ,@('69-43=', '62-39=')\n    ,@('73+21=', '13+55=')\n    ,@('77-49=', '30+0=')\n    ,@('95-83=', '63+32=')\n    ,@('18-12=', '1+85=')\n    ,@('60-55=', '39-28=')\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n    if (-not $found) {\n        throw \"Could not find expected text: $oldText\"\n    }\n}\n\n"}
